# Generate Report for Handback
# - Overview sheet: status text "Ready for handoff" -> "Handed back: in sync with en-US"
# - zh-cn / de-de sheets: status text updated the same way, plus the
#   "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime"
#   columns are now populated now that handback has happened.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

# Status text got longer, widen the status columns
$overview.Columns.Item(5).ColumnWidth = 29.14
$overview.Columns.Item(6).ColumnWidth = 29.14

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

# Status column
$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"

# Row 2 - 8d24fe72-5bf9-46f0-90bf-dbb5e9895d8d
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c14c5649837278297520b161fbcbd2b187a89a1b/e2e/8d24fe72-5bf9-46f0-90bf-dbb5e9895d8d.md", "", "", "8d24fe72-5bf9-46f0-90bf-dbb5e9895d8d.md") | Out-Null
$zhcn.Range("J2").Value = "8d24fe72-5bf9-46f0-90bf-dbb5e9895d8d.c1500da2d1a6a739b0f6fefba8019e4ba5646042.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-09-06 05:21:40"

# Row 3 - d19a2859-f3bf-46b5-97d4-0466a0e4744c
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c14c5649837278297520b161fbcbd2b187a89a1b/e2e/d19a2859-f3bf-46b5-97d4-0466a0e4744c.md", "", "", "d19a2859-f3bf-46b5-97d4-0466a0e4744c.md") | Out-Null
$zhcn.Range("J3").Value = "d19a2859-f3bf-46b5-97d4-0466a0e4744c.3c1b9ff64d8b8abc499c69e0a9b43f81a4b532f5.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-09-06 05:21:40"

# Column widths (status column + the two newly-filled columns need more room)
$zhcn.Columns.Item(3).ColumnWidth = 29.14
$zhcn.Columns.Item(9).ColumnWidth = 39.14
$zhcn.Columns.Item(10).ColumnWidth = 39.14

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

# Status column
$dede.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("C3").Value = "Handed back: in sync with en-US"

# Row 2 - 8d24fe72-5bf9-46f0-90bf-dbb5e9895d8d
$dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c14c5649837278297520b161fbcbd2b187a89a1b/e2e/8d24fe72-5bf9-46f0-90bf-dbb5e9895d8d.md", "", "", "8d24fe72-5bf9-46f0-90bf-dbb5e9895d8d.md") | Out-Null
$dede.Range("J2").Value = "8d24fe72-5bf9-46f0-90bf-dbb5e9895d8d.c1500da2d1a6a739b0f6fefba8019e4ba5646042.de-de.xlf"
$dede.Range("K2").Value = "2016-09-06 05:21:58"

# Row 3 - d19a2859-f3bf-46b5-97d4-0466a0e4744c
$dede.Hyperlinks.Add($dede.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c14c5649837278297520b161fbcbd2b187a89a1b/e2e/d19a2859-f3bf-46b5-97d4-0466a0e4744c.md", "", "", "d19a2859-f3bf-46b5-97d4-0466a0e4744c.md") | Out-Null
$dede.Range("J3").Value = "d19a2859-f3bf-46b5-97d4-0466a0e4744c.3c1b9ff64d8b8abc499c69e0a9b43f81a4b532f5.de-de.xlf"
$dede.Range("K3").Value = "2016-09-06 05:21:58"

# Column widths (status column + the two newly-filled columns need more room)
$dede.Columns.Item(3).ColumnWidth = 29.14
$dede.Columns.Item(9).ColumnWidth = 39.14
$dede.Columns.Item(10).ColumnWidth = 39.14

Write-Output "Handback report generated"
